# "added multiple benchmarks (DFT, SFT)"
#
# Parameters sheet: rows 2-4 get new mean (C) values and the std (D)
# column becomes a formula (=Cn/5) instead of a static value. Two new
# rows (5 "HW", 6 "SW") are populated the same way the existing A/M/BUS
# rows are. The sheet selection moves to D7.
#
# Properties sheet: grows from 3 data rows (1,2,3) to 25 data rows
# (1000..25000, step 1000). The sheet selection moves to A26 with the
# whole A1:A26 range marked as selected.

$wb = $excel.ActiveWorkbook

$params = $wb.Worksheets.Item("Parameters")
$props  = $wb.Worksheets.Item("Properties")

# --- Parameters sheet -------------------------------------------------

# Row 2 ("A") keeps its mean, but std is now a formula.
$params.Range("C2").Value = 0.0001
$params.Range("D2").Formula = "=C2/5"

# Row 3 ("M") gets a new mean, std becomes a formula.
$params.Range("C3").Value = 0.00006
$params.Range("D3").Formula = "=C3/5"

# Row 4 ("BUS") gets a new mean, std becomes a formula.
$params.Range("C4").Value = 0.000001
$params.Range("D4").Formula = "=C4/5"

# Row 5 - new "HW" benchmark entry.
$params.Range("A5").Value = "HW"
$params.Range("B5").Value = "gaussian"
$params.Range("C5").Value = 0.00005
$params.Range("D5").Formula = "=C5/5"
$params.Range("E5").Value = $false

# Row 6 - new "SW" benchmark entry.
$params.Range("A6").Value = "SW"
$params.Range("B6").Value = "gaussian"
$params.Range("C6").Value = 0.00006
$params.Range("D6").Formula = "=C6/5"
$params.Range("E6").Value = $false

# Rename the old A_11/A_12/A_1S benchmark names to the new short names.
$params.Range("A2").Value = "A"
$params.Range("A3").Value = "M"
$params.Range("A4").Value = "BUS"

# --- Properties sheet ---------------------------------------------------

for ($i = 1; $i -le 25; $i++) {
    $row = $i + 1
    $props.Cells.Item($row, 1).Value = $i * 1000
}

$props.Activate() | Out-Null
$props.Range("A1:A26").Select() | Out-Null

# Parameters is the tab that stays active/selected in the workbook, with
# D7 as the selected cell there.
$params.Activate() | Out-Null
$params.Range("D7").Select() | Out-Null
